$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# A3 already holds "1234567890" as a shared-string text value (no special
# cell style applied). Copying it into the new cells preserves that same
# text representation (t="s") rather than Excel auto-converting a literal
# "1234567890" string into a numeric value.
$ws.Range("A3").Copy($ws.Range("A19"))
$ws.Range("A3").Copy($ws.Range("A20"))
